# Auto-generated edit script applying the Ridill_Profits.xlsx diff
# For each affected row (identified by sheet + row number), update columns H-N
# with the new values from the commit. Cells that are removed in the diff are
# cleared; cells newly introduced in the diff are set for the first time.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 62
$ws.Range("H62").Value = 145842290
$ws.Range("I62").Value = 83345580
$ws.Range("J62").Value = 208339000
$ws.Range("K62").Value = 83345580
$ws.Range("L62").Value = 208339000
$ws.Range("M62").Value = -83344956
$ws.Range("N62").Value = -208340248

# Row 65
$ws.Range("H65").Value = 145842290
$ws.Range("I65").Value = 83345580
$ws.Range("J65").Value = 208339000
$ws.Range("K65").Value = 416727900
$ws.Range("L65").Value = 1041695000
$ws.Range("M65").Value = -416724780
$ws.Range("N65").Value = -1041701240

# Row 92
$ws.Range("H92").Value = 1009.6957
$ws.Range("I92").Value = 1080.1578
$ws.Range("J92").Value = 675
$ws.Range("K92").Value = 1080.1578
$ws.Range("L92").Value = 675
$ws.Range("M92").Value = 167.8422
$ws.Range("N92").Value = -3171

# Row 137
$ws.Range("H137").Value = 11031117
$ws.Range("J137").Value = 15968463
$ws.Range("L137").Value = 47905389
$ws.Range("N137").Value = -47910489

# Row 141
$ws.Range("H141").Value = 1926.6897
$ws.Range("I141").Value = 726.5
$ws.Range("J141").Value = 3890.6365
$ws.Range("K141").Value = 2179.5
$ws.Range("L141").Value = 11671.9095
$ws.Range("M141").Value = 3000.5
$ws.Range("N141").Value = -22031.9095


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

# Row 37
$ws.Range("H37").Value = 12139.818
$ws.Range("J37").Value = 21407.6
$ws.Range("L37").Value = 21407.6
$ws.Range("N37").Value = -21953.6

# Row 61
$ws.Range("H61").Value = 5309090
$ws.Range("I61").Value = 2771572
$ws.Range("J61").Value = 29415514
$ws.Range("K61").Value = 2771572
$ws.Range("L61").Value = 29415514
$ws.Range("M61").Value = -2771360
$ws.Range("N61").Value = -29415938

# Row 74
$ws.Range("H74").Value = 6454465.5
$ws.Range("I74").Value = 1383.5454
$ws.Range("J74").Value = 22228666
$ws.Range("K74").Value = 1383.5454
$ws.Range("L74").Value = 22228666
$ws.Range("M74").Value = -509.5454
$ws.Range("N74").Value = -22230414

# Row 77
$ws.Range("H77").Value = 6454465.5
$ws.Range("I77").Value = 1383.5454
$ws.Range("J77").Value = 22228666
$ws.Range("K77").Value = 6917.727
$ws.Range("L77").Value = 111143330
$ws.Range("M77").Value = -2549.727
$ws.Range("N77").Value = -111152066

# Row 102
$ws.Range("H102").Value = 2672
$ws.Range("I102").Value = 1915
$ws.Range("J102").Value = 3176.6667
$ws.Range("K102").Value = 1915
$ws.Range("L102").Value = 3176.6667
$ws.Range("M102").Value = -293
$ws.Range("N102").Value = -6420.6667

# Row 136
$ws.Range("H136").Value = 5309090
$ws.Range("I136").Value = 2771572
$ws.Range("J136").Value = 29415514
$ws.Range("K136").Value = 8314716
$ws.Range("L136").Value = 88246542
$ws.Range("M136").Value = -8312166
$ws.Range("N136").Value = -88251642


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 15
$ws.Range("H15").Value = 11750
$ws.Range("I15").Value = 11750
$ws.Range("K15").Value = 11750
$ws.Range("M15").Value = -11523


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

# Row 31
$ws.Range("H31").Value = 4712040
$ws.Range("I31").Value = 2778746.5
$ws.Range("J31").Value = 12997584
$ws.Range("K31").Value = 2778746.5
$ws.Range("L31").Value = 12997584
$ws.Range("M31").Value = -2778451.5
$ws.Range("N31").Value = -12998174

# Row 34
$ws.Range("H34").Value = 4712040
$ws.Range("I34").Value = 2778746.5
$ws.Range("J34").Value = 12997584
$ws.Range("K34").Value = 2778746.5
$ws.Range("L34").Value = 12997584
$ws.Range("M34").Value = -2778544.5
$ws.Range("N34").Value = -12997988

# Row 50
$ws.Range("H50").Value = 13998
$ws.Range("J50").Value = 13998
$ws.Range("L50").Value = 13998
$ws.Range("N50").Value = -15248

# Row 51
$ws.Range("H51").Value = 28928.572
$ws.Range("J51").Value = 28928.572
$ws.Range("L51").Value = 28928.572
$ws.Range("N51").Value = -30400.572

# Row 59
$ws.Range("H59").Value = 24000
$ws.Range("J59").Value = 24000
$ws.Range("L59").Value = 24000
$ws.Range("N59").Value = -26290

# Row 60
$ws.Range("H60").Value = 13666.667
$ws.Range("J60").Value = 20000
$ws.Range("L60").Value = 20000
$ws.Range("N60").Value = -21022

# Row 61
$ws.Range("H61").Value = 28928.572
$ws.Range("J61").Value = 28928.572
$ws.Range("L61").Value = 28928.572
$ws.Range("N61").Value = -29624.572

# Row 74
$ws.Range("H74").Value = 33333.332
$ws.Range("J74").Value = 33333.332
$ws.Range("L74").Value = 33333.332
$ws.Range("N74").Value = -35081.332

# Row 77
$ws.Range("H77").Value = 33333.332
$ws.Range("J77").Value = 33333.332
$ws.Range("L77").Value = 99999.99600000001
$ws.Range("N77").Value = -108735.996

# Row 105
$ws.Range("H105").Value = 7291.3125
$ws.Range("I105").Value = 1685
$ws.Range("J105").Value = 16635.166
$ws.Range("K105").Value = 1685
$ws.Range("L105").Value = 16635.166
$ws.Range("M105").Value = 62
$ws.Range("N105").Value = -20129.166


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 113
$ws.Range("H113").Value = 3270.7058
$ws.Range("I113").Value = 1981.625
$ws.Range("K113").Value = 5944.875
$ws.Range("M113").Value = -3774.875

# Row 131
$ws.Range("H131").Value = 33400868
$ws.Range("I131").Value = 62625270
$ws.Range("J131").Value = 1547.7142
$ws.Range("K131").Value = 187875810
$ws.Range("L131").Value = 4643.142599999999
$ws.Range("M131").Value = -187870770
$ws.Range("N131").Value = -14723.1426


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 101
$ws.Range("H101").Value = 39000
$ws.Range("J101").Value = 39000
$ws.Range("L101").Value = 39000
$ws.Range("N101").Value = -45490

# Row 132
$ws.Range("H132").Value = 9958579
$ws.Range("I132").Value = 8845059
$ws.Range("J132").Value = 15155002
$ws.Range("K132").Value = 26535177
$ws.Range("L132").Value = 45465006
$ws.Range("M132").Value = -26532647
$ws.Range("N132").Value = -45470066


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 122
$ws.Range("H122").Value = 3087821.5
$ws.Range("I122").Value = 3705002.5
$ws.Range("J122").Value = 1916.6666
$ws.Range("K122").Value = 11115007.5
$ws.Range("L122").Value = 5749.9998
$ws.Range("M122").Value = -11112557.5
$ws.Range("N122").Value = -10649.9998

# Row 135
$ws.Range("H135").Value = 36002.668
$ws.Range("J135").Value = 36002.668
$ws.Range("L135").Value = 36002.668
$ws.Range("N135").Value = -46142.668

